# Update header labels (E1:L1) on every worksheet to the more descriptive
# "severity level(s)" wording, then drop the now-redundant "Category" column
# (column M), which held a constant "secondary" marker for every data row.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("E1").Value2 = "% severity levels 1-2"
    $ws.Range("F1").Value2 = "# severity levels 1-2"
    $ws.Range("G1").Value2 = "% severity level 3"
    $ws.Range("H1").Value2 = "# severity level 3"
    $ws.Range("I1").Value2 = "% severity level 4"
    $ws.Range("J1").Value2 = "# severity level 4"
    $ws.Range("K1").Value2 = "% severity level 5"
    $ws.Range("L1").Value2 = "# severity level 5"

    $ws.Columns.Item(13).Delete()
}
